# "rebuild files as classes" - the reference data (user names, sites, group
# names, role names) was re-emitted from Python classes whose __str__/repr
# wrapped the value in literal double quotes, so every free-text value in
# the lookup tables gained a leading/trailing " character. Numeric id
# columns and the column headers are untouched.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: UserTable (user_name, site) ---------------------------------
$ws1 = $wb.Worksheets.Item(1)

# A2/A3 are typed with a leading apostrophe so Excel treats the literal
# leading " as a force-text prefix; this is what flips on the quotePrefix
# cell style that the target workbook shows for those two cells.
$ws1.Range("A2").Value = "'""Bob"""
$ws1.Range("A3").Value = "'""Mary"""
$ws1.Range("A4").Value = '"John"'
$ws1.Range("A5").Value = '"Frank"'
$ws1.Range("A6").Value = '"Sally"'

$ws1.Range("B2").Value = '"Lake Mary"'
$ws1.Range("B3").Value = '"Winter Park"'
$ws1.Range("B4").Value = '"Lake Mary"'
$ws1.Range("B5").Value = '"Sanford"'
$ws1.Range("B6").Value = '"Winter Park"'

# --- Sheet 2: GroupTable (group_id, group_name) ---------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = '"TeamExtreme"'
$ws2.Range("B3").Value = '"Lab42"'
$ws2.Range("B4").Value = '"Misfits"'

# --- Sheet 3: GroupJunctionTable (username, group_id) ---------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").Value = '"Bob"'
$ws3.Range("A3").Value = '"Mary"'
$ws3.Range("A4").Value = '"Mary"'
$ws3.Range("A5").Value = '"Mary"'
$ws3.Range("A6").Value = '"John"'
$ws3.Range("A7").Value = '"Frank"'
$ws3.Range("A8").Value = '"Frank"'
$ws3.Range("A9").Value = '"Sally"'

# --- Sheet 4: RoleTable (role_id, role_name) ------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = '"chemist"'
$ws4.Range("B3").Value = '"biologist"'
$ws4.Range("B4").Value = '"analyst"'

# --- Sheet 5: RoleJunctionTable (username, role_id) -----------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("A2").Value = '"Bob"'
$ws5.Range("A3").Value = '"Mary"'
$ws5.Range("A4").Value = '"Mary"'
$ws5.Range("A5").Value = '"John"'
$ws5.Range("A6").Value = '"John"'
$ws5.Range("A7").Value = '"Frank"'
$ws5.Range("A8").Value = '"Sally"'
$ws5.Range("A9").Value = '"Sally"'

# --- Restore per-sheet selections, finishing on UserTable as active tab ---
$ws2.Range("B4").Select()
$ws3.Range("D25").Select()
$ws4.Range("B4").Select()
$ws5.Range("D15").Select()
$ws1.Range("C14").Select()
